$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Tue Nov 04 21:26:09 EST 2025"
$ws.Range("B3").Value = "Tue Nov 04 21:26:18 EST 2025"
$ws.Range("B4").Value = "Tue Nov 04 21:26:42 EST 2025"
$ws.Range("B5").Value = "Tue Nov 04 21:26:49 EST 2025"
$ws.Range("B6").Value = "Tue Nov 04 21:26:55 EST 2025"
$ws.Range("B7").Value = "Tue Nov 04 21:27:04 EST 2025"

$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Tue Nov 04 21:27:11 EST 2025"
$ws.Range("B3").Value = "Tue Nov 04 21:27:18 EST 2025"
$ws.Range("B4").Value = "Tue Nov 04 21:27:25 EST 2025"
$ws.Range("B5").Value = "Tue Nov 04 21:27:30 EST 2025"
$ws.Range("B6").Value = "Tue Nov 04 21:27:37 EST 2025"
$ws.Range("B7").Value = "Tue Nov 04 21:27:46 EST 2025"
$ws.Range("B8").Value = "Tue Nov 04 21:27:52 EST 2025"
$ws.Range("B9").Value = "Tue Nov 04 21:27:58 EST 2025"
$ws.Range("B10").Value = "Tue Nov 04 21:28:04 EST 2025"
$ws.Range("B11").Value = "Tue Nov 04 21:28:11 EST 2025"
$ws.Range("B12").Value = "Tue Nov 04 21:28:17 EST 2025"
$ws.Range("B13").Value = "Tue Nov 04 21:28:22 EST 2025"
$ws.Range("B14").Value = "Tue Nov 04 21:28:29 EST 2025"
$ws.Range("B15").Value = "Tue Nov 04 21:28:35 EST 2025"
$ws.Range("B16").Value = "Tue Nov 04 21:28:41 EST 2025"
$ws.Range("B17").Value = "Tue Nov 04 21:28:47 EST 2025"
$ws.Range("B18").Value = "Tue Nov 04 21:28:53 EST 2025"
$ws.Range("B20").Value = "Tue Nov 04 21:28:59 EST 2025"
$ws.Range("B21").Value = "Tue Nov 04 21:29:05 EST 2025"

$ws = $wb.Worksheets.Item("Extension")
$ws.Range("B2").Value = "Tue Nov 04 21:29:11 EST 2025"
$ws.Range("B3").Value = "Tue Nov 04 21:29:17 EST 2025"
$ws.Range("B4").Value = "Tue Nov 04 21:29:22 EST 2025"
$ws.Range("B5").Value = "Tue Nov 04 21:29:29 EST 2025"
$ws.Range("B6").Value = "Tue Nov 04 21:29:34 EST 2025"
$ws.Range("B7").Value = "Tue Nov 04 21:29:40 EST 2025"

$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Tue Nov 04 21:29:47 EST 2025"
$ws.Range("B3").Value = "Tue Nov 04 21:29:52 EST 2025"
$ws.Range("B4").Value = "Tue Nov 04 21:29:58 EST 2025"
$ws.Range("B5").Value = "Tue Nov 04 21:30:04 EST 2025"
$ws.Range("B6").Value = "Tue Nov 04 21:30:10 EST 2025"
$ws.Range("B7").Value = "Tue Nov 04 21:30:15 EST 2025"
$ws.Range("B8").Value = "Tue Nov 04 21:30:22 EST 2025"
$ws.Range("B9").Value = "Tue Nov 04 21:30:28 EST 2025"
$ws.Range("B10").Value = "Tue Nov 04 21:30:34 EST 2025"
$ws.Range("B11").Value = "Tue Nov 04 21:30:40 EST 2025"
$ws.Range("B12").Value = "Tue Nov 04 21:30:46 EST 2025"
$ws.Range("B13").Value = "Tue Nov 04 21:30:52 EST 2025"
$ws.Range("B14").Value = "Tue Nov 04 21:30:58 EST 2025"
$ws.Range("B15").Value = "Tue Nov 04 21:31:04 EST 2025"
$ws.Range("B16").Value = "Tue Nov 04 21:31:10 EST 2025"
$ws.Range("B17").Value = "Tue Nov 04 21:31:16 EST 2025"
$ws.Range("B18").Value = "Tue Nov 04 21:31:22 EST 2025"
$ws.Range("B19").Value = "Tue Nov 04 21:31:28 EST 2025"
$ws.Range("B20").Value = "Tue Nov 04 21:31:34 EST 2025"
$ws.Range("B21").Value = "Tue Nov 04 21:31:41 EST 2025"
$ws.Range("B22").Value = "Tue Nov 04 21:31:47 EST 2025"
$ws.Range("B23").Value = "Tue Nov 04 21:31:53 EST 2025"
$ws.Range("B24").Value = "Tue Nov 04 21:31:59 EST 2025"
$ws.Range("B25").Value = "Tue Nov 04 21:32:05 EST 2025"
$ws.Range("B26").Value = "Tue Nov 04 21:32:11 EST 2025"
$ws.Range("B27").Value = "Tue Nov 04 21:32:17 EST 2025"
$ws.Range("B28").Value = "Tue Nov 04 21:32:23 EST 2025"
$ws.Range("B29").Value = "Tue Nov 04 21:32:29 EST 2025"
$ws.Range("B30").Value = "Tue Nov 04 21:32:35 EST 2025"
$ws.Range("B31").Value = "Tue Nov 04 21:32:41 EST 2025"
$ws.Range("B32").Value = "Tue Nov 04 21:32:46 EST 2025"
$ws.Range("B33").Value = "Tue Nov 04 21:32:52 EST 2025"
$ws.Range("B34").Value = "Tue Nov 04 21:32:58 EST 2025"
$ws.Range("B35").Value = "Tue Nov 04 21:33:04 EST 2025"
$ws.Range("B36").Value = "Tue Nov 04 21:33:10 EST 2025"
$ws.Range("B37").Value = "Tue Nov 04 21:33:16 EST 2025"
$ws.Range("B38").Value = "Tue Nov 04 21:33:22 EST 2025"
$ws.Range("B39").Value = "Tue Nov 04 21:33:32 EST 2025"
$ws.Range("B40").Value = "Tue Nov 04 21:33:39 EST 2025"
$ws.Range("B41").Value = "Tue Nov 04 21:33:46 EST 2025"
$ws.Range("B42").Value = "Tue Nov 04 21:33:52 EST 2025"
$ws.Range("B43").Value = "Tue Nov 04 21:33:59 EST 2025"
$ws.Range("B44").Value = "Tue Nov 04 21:34:05 EST 2025"
$ws.Range("B45").Value = "Tue Nov 04 21:34:10 EST 2025"
$ws.Range("B46").Value = "Tue Nov 04 21:34:17 EST 2025"
$ws.Range("B47").Value = "Tue Nov 04 21:34:23 EST 2025"
$ws.Range("B48").Value = "Tue Nov 04 21:34:29 EST 2025"
$ws.Range("B49").Value = "Tue Nov 04 21:34:35 EST 2025"
$ws.Range("B50").Value = "Tue Nov 04 21:34:41 EST 2025"
$ws.Range("B51").Value = "Tue Nov 04 21:34:47 EST 2025"
$ws.Range("B52").Value = "Tue Nov 04 21:34:52 EST 2025"
$ws.Range("B53").Value = "Tue Nov 04 21:34:59 EST 2025"
$ws.Range("B54").Value = "Tue Nov 04 21:35:05 EST 2025"
$ws.Range("B56").Value = "Tue Nov 04 21:35:10 EST 2025"
$ws.Range("B57").Value = "Tue Nov 04 21:35:21 EST 2025"
$ws.Range("B59").Value = "Tue Nov 04 21:35:28 EST 2025"
$ws.Range("B60").Value = "Tue Nov 04 21:35:34 EST 2025"

$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Tue Nov 04 21:35:47 EST 2025"
$ws.Range("B3").Value = "Tue Nov 04 21:35:53 EST 2025"
$ws.Range("B4").Value = "Tue Nov 04 21:35:59 EST 2025"
$ws.Range("B5").Value = "Tue Nov 04 21:36:05 EST 2025"
$ws.Range("B6").Value = "Tue Nov 04 21:36:11 EST 2025"

$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Tue Nov 04 21:36:16 EST 2025"
$ws.Range("B3").Value = "Tue Nov 04 21:36:30 EST 2025"
$ws.Range("B4").Value = "Tue Nov 04 21:36:45 EST 2025"
$ws.Range("B5").Value = "Tue Nov 04 21:36:59 EST 2025"
$ws.Range("B6").Value = "Tue Nov 04 21:37:14 EST 2025"

$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Tue Nov 04 21:35:41 EST 2025"
